$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Fumée de chandelle</m>" ->
#           "Fumée de " + <tl> + "chandelle" + </tl> + </m>
#   ("<tl>...</tl>" wraps just "chandelle", in the existing blue/Courier-New
#    "tag" run style that the "</m>" run already carries.)
# ---------------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.ClearFormatting()
$f1.Find.Execute("Fumée de chandelle</m>") | Out-Null
$start1 = $f1.Start

# Borrow the already-correctly-formatted blue/Courier-New run ("</m>") as the
# clipboard source for the new "<tl>" tag run, so the new run picks up every
# rPr child (rFonts incl. eastAsia/cs, color, sz, szCs, rtl) exactly.
$srcTl1 = $d.Range($start1 + 18, $start1 + 22)
$srcTl1.Copy()
$dstTl1 = $d.Range($start1 + 9, $start1 + 9)
$dstTl1.Paste()
$pastedTl1 = $d.Range($start1 + 9, $start1 + 13)
$pastedTl1.Text = "<tl>"

# The trailing run (currently "</m>") now sits right after "chandelle"; turn
# it into "</tl></m>" in place so it keeps its original rPr untouched.
$tailRng1 = $d.Range($start1 + 22, $start1 + 26)
$tailRng1.Text = "</tl></m>"

# ---------------------------------------------------------------------------
# Change 2: "moules de noyau</tl> pour" -> "moules</tl> de noyau pour"
#   (Plain text edits inside two already-black runs; no formatting changes.)
# ---------------------------------------------------------------------------
$f2 = $d.Content
$f2.Find.ClearFormatting()
$f2.Find.Execute("moules de noyau</tl> pour") | Out-Null
$start2 = $f2.Start

$r1_2 = $d.Range($start2, $start2 + 15)
$r1_2.Text = "moules"

$start2b = $start2 + 6 + 5
$r2_2 = $d.Range($start2b, $start2b + 5)
$r2_2.Text = " de noyau pour"

# ---------------------------------------------------------------------------
# Change 3: "gect aye la moictié de son espesseur. Mays si elle" ->
#           "gect aye la " + <ms> + "moictié de son espesseur" + </ms> + ". Mays si elle"
# ---------------------------------------------------------------------------
$srcMsOpen = $d.Content
$srcMsOpen.Find.ClearFormatting()
$srcMsOpen.Find.Execute("meslé de <ms>moictié</ms>") | Out-Null
$srcMsOpenStart = $srcMsOpen.Start
$msOpenSrc = $d.Range($srcMsOpenStart + 9, $srcMsOpenStart + 13)

$srcMsClose = $d.Content
$srcMsClose.Find.ClearFormatting()
$srcMsClose.Find.Execute("moictié</ms>") | Out-Null
$srcMsCloseStart = $srcMsClose.Start
$msCloseSrc = $d.Range($srcMsCloseStart + 7, $srcMsCloseStart + 12)

$f3 = $d.Content
$f3.Find.ClearFormatting()
$f3.Find.Execute("gect aye la moictié de son espesseur. Mays si elle") | Out-Null
$start3 = $f3.Start

$msOpenSrc.Copy()
$dstOpen3 = $d.Range($start3 + 12, $start3 + 12)
$dstOpen3.Paste()

$msCloseSrc.Copy()
$dstClose3 = $d.Range($start3 + 40, $start3 + 40)
$dstClose3.Paste()

# ---------------------------------------------------------------------------
# Change 4: "left-bottom" -> "left-" + "middle"
#   ("middle" is a brand-new run whose rPr only carries rtl (no color), matching
#    the minimal-rPr "plain" run style used in a few other spots of the doc.)
# ---------------------------------------------------------------------------
$srcMiddle = $d.Content
$srcMiddle.Find.ClearFormatting()
$srcMiddle.Find.Execute("Encores que") | Out-Null
$srcMiddleStart = $srcMiddle.Start
$middleSrc = $d.Range($srcMiddleStart, $srcMiddleStart + 1)

$f4 = $d.Content
$f4.Find.ClearFormatting()
$f4.Find.Execute("left-bottom") | Out-Null
$start4 = $f4.Start

$leftRng4 = $d.Range($start4, $start4 + 11)
$leftRng4.Text = "left-"

$middleSrc.Copy()
$dst4 = $d.Range($start4 + 5, $start4 + 5)
$dst4.Paste()
$pasted4 = $d.Range($start4 + 5, $start4 + 6)
$pasted4.Text = "middle"
